$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Date") values were re-entered as plain text strings in
# US mm-dd-yyyy format instead of real Excel dates.
$dates = [ordered]@{
  "E2" = "10-18-2021"
  "E3" = "10-18-2021"
  "E4" = "10-18-2021"
  "E5" = "10-18-2020"
  "E6" = "10-18-2021"
  "E7" = "10-18-2021"
  "E8" = "10-18-2021"
}

foreach ($addr in $dates.Keys) {
  $cell = $ws.Range($addr)
  # Leading apostrophe forces Excel to store the value as literal text
  # instead of re-parsing "10-18-2021" back into a date serial.
  $cell.Value = "'" + $dates[$addr]
  # Drop back to the workbook's default (unstyled) cell format now that
  # the cell no longer holds a date value, matching the original sheet's
  # plain text columns.
  $cell.Style = "Normal"
}
